$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date formatting from the row above (A2) onto the new date cell (A4)
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)

# Fill in the new row of data (row 4): date (03/10/2013) + hours worked
$ws.Range("A4").Value = 41550
$ws.Range("B4").Value = 3.5

# Move the active selection to B12, matching the author's final cursor position
$ws.Range("B12").Select()
